# G RC카 BOM.xlsx - "Add files via upload" edit
#
# The motor-driver line item (row 11) changed from a TB6612FNG board to a
# DRV8838 board, with new unit/total prices. Everything else in the sheet
# keeps its original text - only the shared-string table shrinks by one
# entry (TB6612FNG is removed) and grows by two (the new DRV8838 labels),
# which is why the workbook's active selection also shifts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# H11 must be written before D11 so the new shared strings land in the
# same table order as the source edit ( " DRV8838" then "DRV8838" ).
$ws.Range("H11").Value = " DRV8838"

# Clear() (not just setting .Value) drops D11's old explicit cell style,
# matching the author's edit where D11 loses its style index entirely.
$ws.Range("D11").Clear()
$ws.Range("D11").Value = "DRV8838"

$ws.Range("F11").Value = 8400
$ws.Range("G11").Value = 8400

# Move the active selection the way the saved workbook view recorded it.
[void]$ws.Range("G11").Select()
